$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 2 de Septiembre de 2020 a las 11:47"

# Country-name corrections caused by rank swaps (data refresh)
$ws.Range("A47").Value = "Polonia"
$ws.Range("A48").Value = "Japon"
$ws.Range("A119").Value = "Eslovaquia"
$ws.Range("A120").Value = "Mozambique"
$ws.Range("A130").Value = "Eslovenia"
$ws.Range("A131").Value = "Lituania"

# Updated statistics for the new snapshot
$ws.Range("B4").Value = 6258028
$ws.Range("C4").Value = 457
$ws.Range("E4").Value = 2571690
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 188907
$ws.Range("B6").Value = 3773483
$ws.Range("C6").Value = 7375
$ws.Range("D6").Value = 2902096
$ws.Range("E6").Value = 804896
$ws.Range("G6").Value = 31
$ws.Range("H6").Value = 66491
$ws.Range("B26").Value = 180646
$ws.Range("C26").Value = 3075
$ws.Range("D26").Value = 129971
$ws.Range("E26").Value = 43059
$ws.Range("G26").Value = 111
$ws.Range("H26").Value = 7616
$ws.Range("B29").Value = 119627
$ws.Range("C29").Value = 1089
$ws.Range("D29").Value = 96662
$ws.Range("E29").Value = 22002
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = 963
$ws.Range("B47").Value = 68517
$ws.Range("C47").Value = 595
$ws.Range("D47").Value = 47865
$ws.Range("E47").Value = 18574
$ws.Range("G47").Value = 20
$ws.Range("H47").Value = 2078
$ws.Range("B48").Value = 68392
$ws.Range("D48").Value = 57823
$ws.Range("E48").Value = 9273
$ws.Range("H48").Value = 1296
$ws.Range("B71").Value = 27969
$ws.Range("C71").Value = 327
$ws.Range("D71").Value = 23820
$ws.Range("E71").Value = 3415
$ws.Range("B93").Value = 10114
$ws.Range("C93").Value = 10
$ws.Range("D93").Value = 9347
$ws.Range("E93").Value = 508
$ws.Range("B96").Value = 9360
$ws.Range("C96").Value = 6
$ws.Range("D96").Value = 9079
$ws.Range("E96").Value = 153
$ws.Range("B101").Value = 8161
$ws.Range("C101").Value = 19
$ws.Range("E101").Value = 625
$ws.Range("B112").Value = 4831
$ws.Range("C112").Value = 8
$ws.Range("D112").Value = 4401
$ws.Range("E112").Value = 338
$ws.Range("B119").Value = 4042
$ws.Range("C119").Value = 53
$ws.Range("D119").Value = 2523
$ws.Range("E119").Value = 1486
$ws.Range("H119").Value = 33
$ws.Range("B120").Value = 4039
$ws.Range("D120").Value = 2170
$ws.Range("E120").Value = 1846
$ws.Range("H120").Value = 23
$ws.Range("B130").Value = 2979
$ws.Range("C130").Value = 55
$ws.Range("D130").Value = 2359
$ws.Range("E130").Value = 486
$ws.Range("G130").Value = 1
$ws.Range("H130").Value = 134
$ws.Range("B131").Value = 2958
$ws.Range("C131").Value = 29
$ws.Range("D131").Value = 1874
$ws.Range("E131").Value = 998
$ws.Range("H131").Value = 86
$ws.Range("B166").Value = 935
$ws.Range("C166").Value = 16
$ws.Range("E166").Value = 572
